$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# force-formatted as Text first, so Excel keeps them as strings (matching
# the original workbook, where the whole Price column is stored as text,
# e.g. "18.00" / "1.630" must keep trailing zeros, not become 18 / 1.63).
$textCells = @("D5", "D7", "D9", "D10", "D11", "D12", "D13", "D15", "D17", "D18", "D19", "D21", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.825.40'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '1.873.76'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '301.18'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = '0.5343'
$ws.Range('E7').Value = '  +1.79%  '
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('D9').Value = '0.07178'
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('D10').Value = '21.61'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').Value = '0.8877'
$ws.Range('E11').Value = '  -1.99%  '
$ws.Range('D12').Value = '0.08176'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = '93.63'
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.832.60'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').Value = '5.281'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').Value = '14.77'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = '0.000008551'
$ws.Range('E18').Value = '  -1.58%  '
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').Value = '26.858.34'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').Value = '4.984'
$ws.Range('E21').Value = '  -2.67%  '
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('D24').Value = '146.46'
$ws.Range('E24').Value = '  -1.97%  '
$ws.Range('D25').Value = '2.269'
$ws.Range('E25').Value = '  -3.13%  '
$ws.Range('D26').Value = '1.740'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = '18.00'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').Value = '113.92'
$ws.Range('E28').Value = '  -2.07%  '
$ws.Range('D29').Value = '4.723'
$ws.Range('E29').Value = '  -2.40%  '
$ws.Range('D30').Value = '4.607'
$ws.Range('E30').Value = '  -6.02%  '
$ws.Range('D31').Value = '0.09154'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = '0.8049'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').Value = '0.04971'
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('D34').Value = '1.174'
$ws.Range('E34').Value = '  -4.43%  '
$ws.Range('D35').Value = '2.973'
$ws.Range('E35').Value = '  -0.19%  '
$ws.Range('D36').Value = '0.6005'
$ws.Range('E36').Value = '  +4.88%  '
$ws.Range('D37').Value = '3.219'
$ws.Range('E37').Value = '  -4.54%  '
$ws.Range('D38').Value = '2.592'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('E39').Value = '  -2.17%  '
$ws.Range('D40').Value = '1.071'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '8.878'
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '6.541'
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '114.90'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = '0.5112'
$ws.Range('E44').Value = '  +4.37%  '
$ws.Range('D45').Value = '0.1490'
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '1.630'
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.886'
$ws.Range('E48').Value = '  -2.99%  '
$ws.Range('D49').Value = '37.54'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('D50').Value = '0.06045'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('D51').Value = '62.03'
$ws.Range('E51').Value = '  -3.51%  '
